$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "name" -> "x" and "url" -> "y" ---
$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"

# --- Rebuild rows 2+: drop the old stray-formatted empty cells (B2/E2)
#     and replace with new data rows 2-4 ---
[void]$ws.Rows.Item(2).Delete()

$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 5

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "a"

# --- Column widths (characters) ---
$ws.Columns.Item(1).ColumnWidth = 6.2
$ws.Columns.Item(2).ColumnWidth = 6.68
$ws.Columns.Item(3).ColumnWidth = 23.14
$ws.Columns.Item(4).ColumnWidth = 16.9
$ws.Columns.Item(5).ColumnWidth = 19.99

# --- View: zoom + active selection ---
$excel.ActiveWindow.Zoom = 280
[void]$ws.Range("A1").Select()
